# The workbook's weekly data refresh re-sorted the existing data rows
# (rows 2-22, columns A:R) into a new order. No cell content actually
# changed other than the row ordering, so re-apply the same permutation
# here: newRow[i] = oldRow[mapping[i]]  (1-based row offsets within the
# A2:R22 block, i.e. offset 1 == sheet row 2, offset 21 == sheet row 22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:R22")
$oldData = $rng.Value2

$rowCount = 21
$colCount = 18

# mapping[newOffset] = oldOffset (both 1-based, offset 1 = sheet row 2)
$mapping = @{
    1  = 18
    2  = 19
    3  = 10
    4  = 11
    5  = 20
    6  = 8
    7  = 21
    8  = 17
    9  = 13
    10 = 14
    11 = 15
    12 = 16
    13 = 2
    14 = 3
    15 = 4
    16 = 9
    17 = 12
    18 = 1
    19 = 5
    20 = 6
    21 = 7
}

$newData = New-Object 'object[,]' $rowCount, $colCount

for ($newRow = 1; $newRow -le $rowCount; $newRow++) {
    $oldRow = $mapping[$newRow]
    for ($col = 1; $col -le $colCount; $col++) {
        $newData[($newRow - 1), ($col - 1)] = $oldData[$oldRow, $col]
    }
}

$rng.Value2 = $newData
